# Apply crawl update for clothes_2022-07-27.xlsx
#  1) Refresh the crawl timestamp (column O) for every data row.
#  2) Two product rows (15/16) had their id/title/href/label swapped.
#  3) Rows 50-55 got reshuffled (each row now holds the data that used to
#     live one row below it, with row 55 wrapping around to row 50's
#     original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-07-27 08:48:07"
$newTimestamp = "2022-07-27 20:56:57"

$lastRow = $ws.UsedRange.Rows.Count

# Helper: write a value into a cell, forcing text storage (even for
# numeric-looking ids/prices such as "3305484005" or "7.95") by using the
# classic leading-apostrophe "store as text" trick, just like typing the
# value into Excel directly. Blank values become an explicit empty string.
function Set-TextCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
}

# 1) Update the timestamp column (O) for every data row (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Value() -eq $oldTimestamp) {
        Set-TextCell $ws $r 15 $newTimestamp
    }
}

# 2) Swap rows 15 and 16 (id / title / href / productAriaLabel columns;
#    the remaining columns already matched between the two rows).
$row15_A = $ws.Cells.Item(15, 1).Value()
$row15_B = $ws.Cells.Item(15, 2).Value()
$row15_C = $ws.Cells.Item(15, 3).Value()
$row15_N = $ws.Cells.Item(15, 14).Value()

$row16_A = $ws.Cells.Item(16, 1).Value()
$row16_B = $ws.Cells.Item(16, 2).Value()
$row16_C = $ws.Cells.Item(16, 3).Value()
$row16_N = $ws.Cells.Item(16, 14).Value()

Set-TextCell $ws 15 1 $row16_A
Set-TextCell $ws 15 2 $row16_B
Set-TextCell $ws 15 3 $row16_C
Set-TextCell $ws 15 14 $row16_N

Set-TextCell $ws 16 1 $row15_A
Set-TextCell $ws 16 2 $row15_B
Set-TextCell $ws 16 3 $row15_C
Set-TextCell $ws 16 14 $row15_N

# 3) Rows 50-55: each row takes on the data that used to be on the row
#    below it (columns A through N), and row 55 wraps around to take
#    what used to be row 50's data.
$firstShiftRow = 50
$lastShiftRow = 55
$textColumns = 1,2,3,4,7,8,9,10,11,12,13,14   # A,B,C,D,G,H,I,J,K,L,M,N
$numericColumns = 5,6                          # E,F

$originalRows = @{}
for ($r = $firstShiftRow; $r -le $lastShiftRow; $r++) {
    $rowValues = @{}
    foreach ($c in (1..14)) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $originalRows[$r] = $rowValues
}

for ($r = $firstShiftRow; $r -le $lastShiftRow; $r++) {
    if ($r -lt $lastShiftRow) {
        $sourceRow = $r + 1
    } else {
        $sourceRow = $firstShiftRow
    }
    $sourceValues = $originalRows[$sourceRow]

    foreach ($c in $textColumns) {
        Set-TextCell $ws $r $c $sourceValues[$c]
    }
    foreach ($c in $numericColumns) {
        $ws.Cells.Item($r, $c).Value = $sourceValues[$c]
    }
}
